$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1725.4375
$ws.Range("J17").Value = 1725.4375
$ws.Range("L17").Value = 5176.3125
$ws.Range("N17").Value = -5512.3125
$ws.Range("H32").Value = 8614.923000000001
$ws.Range("J32").Value = 8899.5
$ws.Range("L32").Value = 8899.5
$ws.Range("N32").Value = -9551.5
$ws.Range("H55").Value = 635.2857
$ws.Range("I55").Value = 499
$ws.Range("J55").Value = 658
$ws.Range("K55").Value = 499
$ws.Range("L55").Value = 658
$ws.Range("M55").Value = -285
$ws.Range("N55").Value = -1086
$ws.Range("H64").Value = 4574.25
$ws.Range("I64").Value = 4101.3335
$ws.Range("K64").Value = 4101.3335
$ws.Range("M64").Value = -3853.3335
$ws.Range("H67").Value = 4574.25
$ws.Range("I67").Value = 4101.3335
$ws.Range("K67").Value = 4101.3335
$ws.Range("M67").Value = -3243.3335
$ws.Range("H137").Value = 2403.2307
$ws.Range("I137").Value = 2206.7144
$ws.Range("J137").Value = 2632.5
$ws.Range("K137").Value = 6620.1432
$ws.Range("L137").Value = 7897.5
$ws.Range("M137").Value = -4070.1432
$ws.Range("N137").Value = -12997.5
$ws.Range("H138").Value = 3272.4106
$ws.Range("I138").Value = 2284
$ws.Range("J138").Value = 3290.3818
$ws.Range("K138").Value = 6852
$ws.Range("L138").Value = 9871.145400000001
$ws.Range("M138").Value = -1712
$ws.Range("N138").Value = -20151.1454
$ws.Range("H141").Value = 3854.818
$ws.Range("I141").Value = 3761.5557
$ws.Range("J141").Value = 4274.5
$ws.Range("K141").Value = 11284.6671
$ws.Range("L141").Value = 12823.5
$ws.Range("M141").Value = -6104.667099999999
$ws.Range("N141").Value = -23183.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8857.212
$ws.Range("I32").Value = 4854.3
$ws.Range("J32").Value = 48886.332
$ws.Range("K32").Value = 4854.3
$ws.Range("L32").Value = 48886.332
$ws.Range("M32").Value = -4567.3
$ws.Range("N32").Value = -49460.332
$ws.Range("H74").Value = 74014.88
$ws.Range("I74").Value = 65659.64
$ws.Range("K74").Value = 65659.64
$ws.Range("M74").Value = -64785.64
$ws.Range("H77").Value = 74014.88
$ws.Range("I77").Value = 65659.64
$ws.Range("K77").Value = 328298.2
$ws.Range("M77").Value = -323930.2
$ws.Range("H130").Value = 51230.5
$ws.Range("J130").Value = 51230.5
$ws.Range("L130").Value = 51230.5
$ws.Range("N130").Value = -61270.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 49181
$ws.Range("J62").Value = 49181
$ws.Range("L62").Value = 49181
$ws.Range("N62").Value = -50553
$ws.Range("H65").Value = 49181
$ws.Range("J65").Value = 49181
$ws.Range("L65").Value = 147543
$ws.Range("N65").Value = -154407
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").Value = $null
$ws.Range("H88").Value = 69950
$ws.Range("J88").Value = 69950
$ws.Range("L88").Value = 69950
$ws.Range("N88").Value = -70762
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").Value = $null
$ws.Range("H91").Value = 69950
$ws.Range("J91").Value = 69950
$ws.Range("L91").Value = 69950
$ws.Range("N91").Value = -72758
$ws.Range("H134").Value = 3318
$ws.Range("I134").Value = 3318
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9954
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -7419
$ws.Range("N134").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3268.125
$ws.Range("I16").Value = 3329.8
$ws.Range("K16").Value = 3329.8
$ws.Range("M16").Value = -3042.8
$ws.Range("H31").Value = 1994.5667
$ws.Range("I31").Value = 1374.3
$ws.Range("K31").Value = 1374.3
$ws.Range("M31").Value = -1079.3
$ws.Range("H34").Value = 1994.5667
$ws.Range("I34").Value = 1374.3
$ws.Range("K34").Value = 1374.3
$ws.Range("M34").Value = -1172.3
$ws.Range("H50").Value = 29568.363
$ws.Range("J50").Value = 29525.2
$ws.Range("L50").Value = 29525.2
$ws.Range("N50").Value = -30775.2
$ws.Range("H51").Value = 27096.5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 27096.5
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 27096.5
$ws.Range("M51").Value = $null
$ws.Range("N51").Value = -28568.5
$ws.Range("H61").Value = 27096.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 27096.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 27096.5
$ws.Range("M61").Value = $null
$ws.Range("N61").Value = -27792.5
$ws.Range("H62").Value = 8125.375
$ws.Range("I62").Value = 6100.6
$ws.Range("J62").Value = 11500
$ws.Range("K62").Value = 6100.6
$ws.Range("L62").Value = 11500
$ws.Range("M62").Value = -5476.6
$ws.Range("N62").Value = -12748
$ws.Range("H65").Value = 8125.375
$ws.Range("I65").Value = 6100.6
$ws.Range("J65").Value = 11500
$ws.Range("K65").Value = 30503
$ws.Range("L65").Value = 57500
$ws.Range("M65").Value = -27383
$ws.Range("N65").Value = -63740
$ws.Range("H105").Value = 3001.75
$ws.Range("I105").Value = 993
$ws.Range("K105").Value = 993
$ws.Range("M105").Value = 754
$ws.Range("H107").Value = 31192.678
$ws.Range("I107").Value = 37885.08
$ws.Range("K107").Value = 37885.08
$ws.Range("M107").Value = -35965.08
$ws.Range("H113").Value = 3268.125
$ws.Range("I113").Value = 3329.8
$ws.Range("K113").Value = 3329.8
$ws.Range("M113").Value = -1159.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2554
$ws.Range("I68").Value = 2581.5
$ws.Range("J68").Value = 2499
$ws.Range("K68").Value = 7744.5
$ws.Range("L68").Value = 7497
$ws.Range("M68").Value = -6933.5
$ws.Range("N68").Value = -9119
$ws.Range("H71").Value = 2554
$ws.Range("I71").Value = 2581.5
$ws.Range("J71").Value = 2499
$ws.Range("K71").Value = 23233.5
$ws.Range("L71").Value = 22491
$ws.Range("M71").Value = -19177.5
$ws.Range("N71").Value = -30603
$ws.Range("H102").Value = 5436
$ws.Range("I102").Value = 4923.2
$ws.Range("K102").Value = 14769.6
$ws.Range("M102").Value = -12335.6
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 10977.6
$ws.Range("I4").Value = 10875
$ws.Range("J4").Value = 11003.25
$ws.Range("K4").Value = 10875
$ws.Range("L4").Value = 11003.25
$ws.Range("M4").Value = -10763
$ws.Range("N4").Value = -11227.25
$ws.Range("H80").Value = 7547.364
$ws.Range("I80").Value = 3020.1333
$ws.Range("J80").Value = 17248.572
$ws.Range("K80").Value = 3020.1333
$ws.Range("L80").Value = 17248.572
$ws.Range("M80").Value = -2022.1333
$ws.Range("N80").Value = -19244.572
$ws.Range("H83").Value = 7547.364
$ws.Range("I83").Value = 3020.1333
$ws.Range("J83").Value = 17248.572
$ws.Range("K83").Value = 15100.6665
$ws.Range("L83").Value = 86242.86
$ws.Range("M83").Value = -10108.6665
$ws.Range("N83").Value = -96226.86
$ws.Range("H105").Value = 77924
$ws.Range("J105").Value = 77924
$ws.Range("L105").Value = 77924
$ws.Range("N105").Value = -84912
$ws.Range("H109").Value = 29998
$ws.Range("J109").Value = 29998
$ws.Range("L109").Value = 29998
$ws.Range("N109").Value = -32078
$ws.Range("H128").Value = 45999.2
$ws.Range("J128").Value = 45999.2
$ws.Range("L128").Value = 45999.2
$ws.Range("N128").Value = -55959.2
$ws.Range("H132").Value = 8287.333000000001
$ws.Range("I132").Value = 9117.583000000001
$ws.Range("J132").Value = 4966.3335
$ws.Range("K132").Value = 27352.749
$ws.Range("L132").Value = 14899.0005
$ws.Range("M132").Value = -24822.749
$ws.Range("N132").Value = -19959.0005
$ws.Range("H136").Value = 16520.188
$ws.Range("J136").Value = 16520.188
$ws.Range("L136").Value = 49560.564
$ws.Range("N136").Value = -54660.564
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2002.138
$ws.Range("J46").Value = 2403.2727
$ws.Range("L46").Value = 2403.2727
$ws.Range("N46").Value = -2779.2727
$ws.Range("I68").Value = 1273.5
$ws.Range("J68").Value = 6750
$ws.Range("K68").Value = 1273.5
$ws.Range("L68").Value = 6750
$ws.Range("M68").Value = -524.5
$ws.Range("N68").Value = -8248
$ws.Range("I71").Value = 1273.5
$ws.Range("J71").Value = 6750
$ws.Range("K71").Value = 6367.5
$ws.Range("L71").Value = 33750
$ws.Range("M71").Value = -2623.5
$ws.Range("N71").Value = -41238
$ws.Range("H124").Value = 68300
$ws.Range("J124").Value = 68300
$ws.Range("L124").Value = 68300
$ws.Range("N124").Value = -78120
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1773.9375
$ws.Range("I113").Value = 1669.5834
$ws.Range("J113").Value = 2087
$ws.Range("K113").Value = 5008.7502
$ws.Range("L113").Value = 6261
$ws.Range("M113").Value = -2838.7502
$ws.Range("N113").Value = -10601
$ws.Range("H124").Value = 56676.332
$ws.Range("J124").Value = 56676.332
$ws.Range("L124").Value = 56676.332
$ws.Range("N124").Value = -66496.33199999999
$ws.Range("H140").Value = 123980
$ws.Range("J140").Value = 123980
$ws.Range("L140").Value = 123980
$ws.Range("N140").Value = -134340
